$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B, C, D, E, G). F is unchanged.
$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    3 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    4 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    5 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    6 = @(1.505614041169197, 87981.0709163148,  16.98373111632243, 6.48142807727062,  88006.04168954957)
    7 = @(3.182878228561681, 9.226618575922256, 16.98373111632243, 6.48142807727062,  35.87465599807698)
    8 = @(0.3464964993005633, 1.65323645889881, 16.98373111632243, 6.48142807727062, 25.46489215179242)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G
}

$wb.Save()
